# TaskPlan.xlsx update:
#  - Fill in actual time spent (column J) for the remaining tasks
#  - Add an "Actual Time" column header and an "Actual Project Time (in h)"
#    summary row, matching the pattern already used for the estimated time
#  - Round the existing average project time calculation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the summary row that used to hold "Actual Case" text so it
# becomes the new "Actual Project Time (in h)" summary line (this reuses
# the shared string slot formerly used by J1's old label).
$ws.Range("A22").Value = "Actual Project Time (in h)"

# Rename the J column header from "Actual Case" to "Actual Time".
$ws.Range("J1").Value = "Actual Time"

# Fill in the actual-time values for the remaining tasks (D7-D11).
$ws.Range("J8").Value = 100
$ws.Range("J9").Value = 30
$ws.Range("J10").Value = 25
$ws.Range("J11").Value = 25
$ws.Range("J12").Value = 70

# Round the average project time estimate to 2 decimal places.
$ws.Range("H20").Formula = "=ROUND(MEDIAN(H17, H18)/60, 2)"

# New summary: total actual time logged, converted to hours and rounded.
$ws.Range("J22").Formula = "=ROUND(SUM(J2:J12)/60, 2)"

$ws.Range("G17").Select() | Out-Null

$wb.Save()
